$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for each record (rows 2-29).
# Update each of these date values from 45556 (2024-09-21) to 45557 (2024-09-22).
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45557
}
